# Add duplicate transaction check and enhance M-PESA message parsing
# -> Appends the newly-parsed M-PESA transaction as row 6 of the transactions sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new row's cells to be stored as text (matching the existing rows,
# which are all text values) so values like "1.00", "5/9/25" or "0.00" are not
# auto-converted to numbers/dates by Excel's type inference.
$ws.Range("A6:K6").NumberFormat = "@"

$ws.Range("A6").Value = "TI59KUQVON"
$ws.Range("B6").Value = "1.00"
$ws.Range("C6").Value = "Send Money"
$ws.Range("D6").Value = "FRANKLINE  ATUTI 0794492538"
$ws.Range("E6").Value = "5/9/25"
$ws.Range("F6").Value = "8:46 PM"
$ws.Range("G6").Value = "2505.09"
$ws.Range("H6").Value = "0.00"
$ws.Range("I6").Value = "497978.00"
$ws.Range("J6").Value = "TI59KUQVON Confirmed. Ksh1.00 sent to FRANKLINE  ATUTI 0794492538 on 5/9/25 at 8:46 PM. New M-PESA balance is Ksh2,505.09. Transaction cost, Ksh0.00.  Amount you can transact within the day is 497,978.00. Earn interest daily on Ziidi MMF,Dial *334#, date=1757094378914"
$ws.Range("K6").Value = "2025-09-05 20:47:47"
